{"js": "// Update the date line and the 25 division problems (5 blocks of 5,\n// laid out as a 20-row x 5-col table where every 4th row holds data).\nconst dateResults = context.document.body.search(\"2025-06-28 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\n\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2025-06-29 Sunday\", \"Replace\");\n}\n\nconst table = tables.items[0];\n\n// New values, in document (reading) order: row-block 0 (table rows 0),\n// row-block 1 (table row 4), ... each block has 5 columns.\nconst newValues = [\n  \"622\u00f77=\", \"577\u00f75=\", \"300\u00f75=\", \"255\u00f77=\", \"237\u00f72=\",\n  \"786\u00f77=\", \"599\u00f78=\", \"159\u00f79=\", \"311\u00f78=\", \"497\u00f79=\",\n  \"852\u00f73=\", \"696\u00f79=\", \"688\u00f73=\", \"454\u00f77=\", \"978\u00f78=\",\n  \"787\u00f78=\", \"434\u00f79=\", \"431\u00f79=\", \"750\u00f73=\", \"297\u00f78=\",\n  \"159\u00f77=\", \"768\u00f78=\", \"548\u00f79=\", \"174\u00f72=\", \"289\u00f79=\"\n];\n\nconst dataRows = [0, 4, 8, 12, 16];\nlet i = 0;\nfor (const row of dataRows) {\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems (5 blocks of 5,\n# laid out as a 20-row x 5-col table where every 4th row holds data).\n$d = $word.ActiveDocument\n\n# Date heading (first paragraph).\n$d.Paragraphs.Item(1).Range.Text = \"2025-06-29 Sunday\"\n\n$t = $d.Tables.Item(1)\n\n# New values, in document (reading) order: block 1 = table row 1,\n# block 2 = table row 5, block 3 = row 9, block 4 = row 13, block 5 = row 17.\n$newValues = @(\n    \"622\u00f77=\", \"577\u00f75=\", \"300\u00f75=\", \"255\u00f77=\", \"237\u00f72=\",\n    \"786\u00f77=\", \"599\u00f78=\", \"159\u00f79=\", \"311\u00f78=\", \"497\u00f79=\",\n    \"852\u00f73=\", \"696\u00f79=\", \"688\u00f73=\", \"454\u00f77=\", \"978\u00f78=\",\n    \"787\u00f78=\", \"434\u00f79=\", \"431\u00f79=\", \"750\u00f73=\", \"297\u00f78=\",\n    \"159\u00f77=\", \"768\u00f78=\", \"548\u00f79=\", \"174\u00f72=\", \"289\u00f79=\"\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$i = 0\nforeach ($row in $dataRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($row, $col).Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
